# Apply crypto price/volume table update per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force plain-text storage for numeric-looking price strings so formatting (e.g. trailing zeros) is preserved
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "67.039.23"
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").Value = "3.612.80"
$ws.Range("E3").Value = "  -0.87%  "

$ws.Range("D5").Value = "587.47"
$ws.Range("E5").Value = "  -0.01%  "

$ws.Range("D6").Value = "183.18"
$ws.Range("E6").Value = "  +2.64%  "

$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  -2.45%  "

$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").Value = "0.672"
$ws.Range("E9").Value = "  -5.01%  "

$ws.Range("D10").Value = "53.81"
$ws.Range("E10").Value = "  -2.57%  "

$ws.Range("D11").Value = "0.143"
$ws.Range("E11").Value = "  -9.96%  "

$ws.Range("D12").Value = "0.0000252"
$ws.Range("E12").Value = "  -12.51%  "

$ws.Range("D13").Value = "9.90"
$ws.Range("E13").Value = "  -6.18%  "

$ws.Range("D14").Value = "4.198.47"
$ws.Range("E14").Value = "  -0.15%  "

$ws.Range("D15").Value = "3.617.75"
$ws.Range("E15").Value = "  -0.52%  "

$ws.Range("D16").Value = "0.125"
$ws.Range("E16").Value = "  -0.20%  "

$ws.Range("D17").Value = "67.031.24"
$ws.Range("E17").Value = "  -1.26%  "

$ws.Range("D18").Value = "18.33"
$ws.Range("E18").Value = "  -4.44%  "

$ws.Range("D19").Value = "12.18"
$ws.Range("E19").Value = "  -3.81%  "

$ws.Range("E20").Value = "  -4.78%  "

$ws.Range("D21").Value = "391.98"
$ws.Range("E21").Value = "  -3.84%  "

$ws.Range("D22").Value = "4.30"
$ws.Range("E22").Value = "  -4.88%  "

$ws.Range("D23").Value = "84.77"
$ws.Range("E23").Value = "  -3.61%  "

$ws.Range("E24").Value = "  -4.42%  "

$ws.Range("D25").Value = "12.24"
$ws.Range("E25").Value = "  -3.06%  "

$ws.Range("D26").Value = "6.05"
$ws.Range("E26").Value = "  +0.34%  "

$ws.Range("D27").Value = "10.28"
$ws.Range("E27").Value = "  -3.56%  "

$ws.Range("D28").Value = "3.59"
$ws.Range("E28").Value = "  -10.97%  "

$ws.Range("D29").Value = "8.94"
$ws.Range("E29").Value = "  -4.88%  "

$ws.Range("D30").Value = "31.07"
$ws.Range("E30").Value = "  -4.11%  "

$ws.Range("D31").Value = "6.77"
$ws.Range("E31").Value = "  -5.70%  "

$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").Value = "11.89"
$ws.Range("E32").Value = "  -3.02%  "

$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").Value = "65.14"
$ws.Range("E33").Value = "  +1.15%  "

$ws.Range("D34").Value = "597.15"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("E35").Value = "  -3.36%  "

$ws.Range("D36").Value = "41.31"
$ws.Range("E36").Value = "  -3.27%  "

$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("E38").Value = "  +0.38%  "

$ws.Range("E39").Value = "  -5.57%  "

$ws.Range("D40").Value = "0.0₃0735"
$ws.Range("E40").Value = "  -16.54%  "

$ws.Range("E41").Value = "  -2.44%  "

$ws.Range("E42").Value = "  -7.54%  "

$ws.Range("E43").Value = "  -5.43%  "

$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.733.51"
$ws.Range("E44").Value = "  +2.00%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").Value = "2.41"
$ws.Range("E45").Value = "  -10.35%  "

$ws.Range("E46").Value = "  -3.37%  "

$ws.Range("D47").Value = "3.06"
$ws.Range("E47").Value = "  -1.71%  "

$ws.Range("E48").Value = "  -5.49%  "

$ws.Range("D49").Value = "135.75"
$ws.Range("E49").Value = "  -2.78%  "

$ws.Range("E50").Value = "  -7.64%  "

$ws.Range("E51").Value = "  -6.12%  "
